$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "TLS"
$ws.Range("B6").Value = -8.874217
$ws.Range("C6").Value = 125.727539
$ws.Range("D6").Value = "Timor-Leste"

$ws.Range("A6:D6").Select()
